{"js": "// Locate the paragraph that discusses singly/doubly linked lists by its\n// distinctive trailing text, then grab the (currently empty) paragraph that\n// immediately follows it. Together these two paragraphs are replaced with:\n//   1) the same \"\ub2e8\uc77c/\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8\" paragraph, but without the eastAsia\n//      paragraph-mark hint and without the mid-paragraph \"_GoBack\" bookmark,\n//   2) a new empty paragraph,\n//   3) a new \"\ub2e8\uc77c \uc5f0\uacb0 \ub9ac\uc2a4\ud2b8\" heading paragraph,\n//   4) a new paragraph describing a singly linked list, ending with the\n//      \"_GoBack\" bookmark (moved here from paragraph 1).\nconst body = context.document.body;\n\nconst searchResults = body.search(\"\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target paragraph (\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8 text not found).\");\n}\n\nconst hitParagraphs = searchResults.items[0].paragraphs;\nhitParagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = hitParagraphs.items[0];\nconst followingParagraph = targetParagraph.getNextOrNullObject();\nfollowingParagraph.load(\"isNullObject,text\");\nawait context.sync();\n\nconst startRange = targetParagraph.getRange(\"Whole\");\nconst endRange = (!followingParagraph.isNullObject && followingParagraph.text === \"\")\n  ? followingParagraph.getRange(\"Whole\")\n  : startRange;\nconst replaceRange = startRange.expandTo(endRange);\n\nconst W_NS = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\";\n\nconst bodyFragment =\n  '<w:p xmlns:w=\"' + W_NS + '\">' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uac01 \ub178\ub4dc\uac00 \ub370\uc774\ud130\uc640 \ud3ec\uc778\ud130\ub97c \uac00\uc9c0\uace0 \uc788\uc73c\uba74\uc11c \ub178\ub4dc\ub4e4\uc774 \ud55c \uc904\ub85c \uc5f0\uacb0\ub418\uc5b4 \uc788\ub294 \ubc29\uc2dd.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\ub178\ub4dc\ub4e4\uc774 \ud55c \ubc29\ud5a5\uc73c\ub85c \ub2e4\uc74c \ub178\ub4dc\ub97c \uac00\ub9ac\ud0a4\ub294 \ub9ac\uc2a4\ud2b8\ub97c \ub2e8\uc77c \uc5f0\uacb0 \ub9ac\uc2a4\ud2b8</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc774\uace0</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>,</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uac01 \ub178\ub4dc\uac00 \uc774\uc804 \ub178\ub4dc \ub2e4\uc74c \ub178\ub450\ub97c \ubaa8\ub450 \uac00\ub9ac\ud0a4\ub294 \uc591\ubc29\ud5a5 \ub178\ub4dc\ub85c \ub418\uc5b4 \uc788\ub294 \uac83</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc740</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> \uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p xmlns:w=\"' + W_NS + '\"/>' +\n  '<w:p xmlns:w=\"' + W_NS + '\">' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\ub2e8\uc77c \uc5f0\uacb0 \ub9ac\uc2a4\ud2b8</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p xmlns:w=\"' + W_NS + '\">' +\n    '<w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:lastRenderedPageBreak/><w:t>\ub2e8\ubc29\ud5a5\uc73c\ub85c \ub178\ub4dc\ub4e4\uc744 \uc5f0\uacb0.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\ub178\ub4dc\uc758 \ub370\uc774\ud130\ud53c\ub97b\uc640 \ub2e4\uc74c \ub178\ub4dc\ub97c \uac00\ub9ac\ud0a4\ub294 \ud3ec\uc778\ud130\ub85c \uad6c\uc131.</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> \uccab \ub178\ub4dc\ub97c \uac00\ub9ac\ud0a4\ub294 \ud5e4\ub4dc \ud544\ub4dc\ub97c \uc0ac\uc6a9\ud558\uc5ec \uc804\uccb4 \ub9ac\uc2a4\ud2b8\ub97c \uc21c\ucc28\uc801\uc73c\ub85c \uc5d1\uc138\uc2a4.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>';\n\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n      '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n      '</Relationships></pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData><w:document xmlns:w=\"' + W_NS + '\"><w:body>' + bodyFragment + '</w:body></w:document></pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nreplaceRange.insertOoxml(flatOpcPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Locate the paragraph that discusses singly/doubly linked lists by its\n# distinctive trailing text (\"\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8\"), then find the (currently\n# empty) paragraph that immediately follows it. Together these two\n# paragraphs are replaced with:\n#   1) the same \"\ub2e8\uc77c/\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8\" paragraph, but without the eastAsia\n#      paragraph-mark hint and without the mid-paragraph \"_GoBack\" bookmark,\n#   2) a new empty paragraph,\n#   3) a new \"\ub2e8\uc77c \uc5f0\uacb0 \ub9ac\uc2a4\ud2b8\" heading paragraph,\n#   4) a new paragraph describing a singly linked list, ending with the\n#      \"_GoBack\" bookmark (moved here from paragraph 1).\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8\")\nif (-not $found) {\n    throw \"Could not find the target paragraph (\uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8 text not found).\"\n}\n$hitStart = $searchRange.Start\n\n$paraCount = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $hitStart -and $hitStart -lt $p.Range.End) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not map the found text back to a paragraph.\"\n}\n\n$endIndex = $targetIndex\nif ($targetIndex -lt $paraCount) {\n    $pNext = $d.Paragraphs.Item($targetIndex + 1)\n    if ($pNext.Range.Text.Trim() -eq \"\") {\n        $endIndex = $targetIndex + 1\n    }\n}\n\n$pStart = $d.Paragraphs.Item($targetIndex)\n$pEnd = $d.Paragraphs.Item($endIndex)\n$replaceRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)\n\n$W_NS = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"\n\n$bodyFragment = (\n  '<w:p xmlns:w=\"' + $W_NS + '\">' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uac01 \ub178\ub4dc\uac00 \ub370\uc774\ud130\uc640 \ud3ec\uc778\ud130\ub97c \uac00\uc9c0\uace0 \uc788\uc73c\uba74\uc11c \ub178\ub4dc\ub4e4\uc774 \ud55c \uc904\ub85c \uc5f0\uacb0\ub418\uc5b4 \uc788\ub294 \ubc29\uc2dd.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\ub178\ub4dc\ub4e4\uc774 \ud55c \ubc29\ud5a5\uc73c\ub85c \ub2e4\uc74c \ub178\ub4dc\ub97c \uac00\ub9ac\ud0a4\ub294 \ub9ac\uc2a4\ud2b8\ub97c \ub2e8\uc77c \uc5f0\uacb0 \ub9ac\uc2a4\ud2b8</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc774\uace0</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>,</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uac01 \ub178\ub4dc\uac00 \uc774\uc804 \ub178\ub4dc \ub2e4\uc74c \ub178\ub450\ub97c \ubaa8\ub450 \uac00\ub9ac\ud0a4\ub294 \uc591\ubc29\ud5a5 \ub178\ub4dc\ub85c \ub418\uc5b4 \uc788\ub294 \uac83</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc740</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> \uc774\uc911 \uc5f0\uacb0\ub9ac\uc2a4\ud2b8</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p xmlns:w=\"' + $W_NS + '\"/>' +\n  '<w:p xmlns:w=\"' + $W_NS + '\">' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\ub2e8\uc77c \uc5f0\uacb0 \ub9ac\uc2a4\ud2b8</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p xmlns:w=\"' + $W_NS + '\">' +\n    '<w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:lastRenderedPageBreak/><w:t>\ub2e8\ubc29\ud5a5\uc73c\ub85c \ub178\ub4dc\ub4e4\uc744 \uc5f0\uacb0.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\ub178\ub4dc\uc758 \ub370\uc774\ud130\ud53c\ub97b\uc640 \ub2e4\uc74c \ub178\ub4dc\ub97c \uac00\ub9ac\ud0a4\ub294 \ud3ec\uc778\ud130\ub85c \uad6c\uc131.</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> \uccab \ub178\ub4dc\ub97c \uac00\ub9ac\ud0a4\ub294 \ud5e4\ub4dc \ud544\ub4dc\ub97c \uc0ac\uc6a9\ud558\uc5ec \uc804\uccb4 \ub9ac\uc2a4\ud2b8\ub97c \uc21c\ucc28\uc801\uc73c\ub85c \uc5d1\uc138\uc2a4.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>'\n)\n\n$xml = (\n  '<?xml version=\"1.0\" encoding=\"UTF-8\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData><w:document xmlns:w=\"' + $W_NS + '\"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n)\n\n$replaceRange.InsertXML($xml)\n"}
